# Add a "Summe" (total) row under the "Frage 1" table, summing column C,
# and leave the selection on each sheet where the author last clicked
# (B45 on "Frage 1", C8 on "Frage 2") -- the rest of the commit is just
# Excel-version save noise ("open a few excels and probably just changed
# the selected cell :S").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Frage 1")
$ws1.Range("B44").Value = "Summe"
$ws1.Range("C44").Formula = "=SUM(C5:C41)"

$ws2 = $wb.Worksheets.Item("Frage 2")
[void]$ws2.Range("C8").Select()

# Re-select "Frage 1" last so it stays the active/tab-selected sheet,
# with the cursor parked on B45 as in the saved file.
[void]$ws1.Select()
[void]$ws1.Range("B45").Select()
